# This script updates the "1401ME16" marksheet quiz worksheet so that it
# reflects the actual (previously all-zero / "Absent") results of the
# student, including handling non-integer/float-derived inputs correctly
# (per the commit message "Handles float input without breaking stuff").
#
# Summary of the underlying model (reverse engineered from the target state):
#   - Columns A/B (rows 16-40) and D/E (rows 16-18) each represent one
#     question: column A/D = "Student Ans" (style correctStyle if right,
#     incorrectStyle if wrong, normalStyle/blank if not attempted), column
#     B/E = "Correct Ans" (style absoluteStyle, always shown).
#   - The third block (columns G/H) is no longer used and is removed.
#   - Row 10 = counts of Right/Wrong/Not-Attempted/Max.
#   - Row 11 = marks awarded per right/wrong answer.
#   - Row 12 = total marks for right/wrong answers and the final score
#     string "earned/max".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 10: Right / Wrong / Not Attempt / Max counts
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "No."
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 13
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

# ---------------------------------------------------------------------
# Row 11: marking scheme (marks per right / wrong answer)
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Marking"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# ---------------------------------------------------------------------
# Row 12: totals and final score
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "Total"
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -5
$ws.Range("E12").Value = "47/112"

# ---------------------------------------------------------------------
# Third "Student Ans / Correct Ans" block (columns G/H) is no longer
# used - remove it entirely (rows 15-40).
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------
# Second block (columns D/E): only the first three questions remain
# (rows 16-18); the rest (rows 19-40) are removed.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

# Update the student answers that are now known for the remaining D/E
# rows (row 17 stayed "not attempted").
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

# ---------------------------------------------------------------------
# First block (columns A/B): fill in the student's actual answers for
# each of the 25 questions (rows 16-40), with the appropriate style:
#   correctStyle   -> answered correctly
#   incorrectStyle -> answered incorrectly
#   normalStyle    -> left blank (not attempted) - already the default
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"

$ws.Range("A17").Value = "Option D"
$ws.Range("A17").Style = "correctStyle"

$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"

$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"

$ws.Range("A22").Value = "Option A"
$ws.Range("A22").Style = "incorrectStyle"

$ws.Range("A25").Value = "Option B"
$ws.Range("A25").Style = "incorrectStyle"

$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"

$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"

$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"

$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"

$ws.Range("A34").Value = "Option A"
$ws.Range("A34").Style = "incorrectStyle"

$ws.Range("A35").Value = "Option D"
$ws.Range("A35").Style = "correctStyle"

$ws.Range("A36").Value = "Option A"
$ws.Range("A36").Style = "correctStyle"

$ws.Range("A37").Value = "Option B"
$ws.Range("A37").Style = "incorrectStyle"

$ws.Range("A38").Value = "Option B"
$ws.Range("A38").Style = "incorrectStyle"

$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
